$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 143, shifting rows 143:258 down to 144:259
$ws.Rows.Item(143).Insert()

# Populate the new row 143 with the data (copied structure from neighboring rows,
# with date/value/origin fields updated per the target change)
$ws.Range("A143").Value = 10
$ws.Range("B143").Value = "Vega Modelo de Temuco"
$ws.Range("C143").Value = "La Araucanía"
$ws.Range("D143").Value = 44484
$ws.Range("D143").NumberFormat = $ws.Range("D144").NumberFormat
$ws.Range("E143").Value = 9
$ws.Range("F143").Value = 100114014
$ws.Range("G143").Value = "Betarraga"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 50
$ws.Range("K143").Value = 9500
$ws.Range("L143").Value = 9500
$ws.Range("M143").Value = 9500
$ws.Range("N143").Value = "$/docena de paquetes"
$ws.Range("O143").Value = "Región del Maule"
$ws.Range("P143").Value = 792
$ws.Range("Q143").Value = 12
$ws.Range("R143").Value = "Hortaliza"
